$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 5 (CreateDate) - change 形態(type) from DATE to TIMESTAMP and clear 長度(length)
$ws.Range("D13").Value = "TIMESTAMP"
$ws.Range("E13").ClearContents()

# Row 7 (LastUpdate) - change 形態(type) from DATE to TIMESTAMP and clear 長度(length)
$ws.Range("D15").Value = "TIMESTAMP"
$ws.Range("E15").ClearContents()

# Move the active selection to B9, matching the author's final cursor position
$ws.Range("B9").Select()
